$d = $word.ActiveDocument

# --------------------------------------------------------------------
# The visible content edit made by the author: "预警" -> "限流" in the
# infrastructure paragraph ("...心跳、监控、预警等方式维系其稳定").
#
# Real Word also keeps a single "_GoBack" bookmark that always marks the
# location of the most recent edit; it previously sat between "团队" and
# "绩效" (a leftover from an earlier edit) and must move to the new edit
# site, collapsing the old run split in the process.
# --------------------------------------------------------------------

# Step 1: re-write "团队绩效" (identical text) through Find/Replace so the
# stale "_GoBack" bookmark sitting inside it is dropped and the two runs
# it used to separate ("团队" / "绩效") are rebuilt as a single run.
$d.Content.Find.Execute("团队绩效", $true, $false, $false, $false, $false,
                         $true, 1, $false, "团队绩效", 2)

# Step 2: the actual text change.
$find = $d.Content.Find
$find.Execute("预警", $true, $false, $false, $false, $false,
              $true, 1, $false, "", 0)
$editRange = $find.Parent
$editRange.Text = "限流"
$editEnd = $editRange.End

# Step 3: drop a fresh "_GoBack" bookmark right after the new text, which
# is where Word leaves it after typing -- this also splits the run at
# that point, matching Word's own behaviour.
$goBack = $d.Range($editEnd, $editEnd)
$d.Bookmarks.Add("_GoBack", $goBack)
